$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for 28fbfdb8-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-20 06:51:32"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 28fbfdb8-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-20 06:51:27"
$wsZhCn.Range("K3").Value = "2016-08-20 06:51:43"

# de-de sheet: Correspond Handback DateTime for 28fbfdb8-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-20 06:51:49"
